$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 153.75
$ws.Range("I5").Value = 121.666664
$ws.Range("J5").Value = 250
$ws.Range("K5").Value = 121.666664
$ws.Range("L5").Value = 250
$ws.Range("M5").Value = -6.666663999999997
$ws.Range("N5").Value = -480

$ws.Range("H12").Value = 466.5625
$ws.Range("J12").Value = 1025
$ws.Range("L12").Value = 1025
$ws.Range("N12").Value = -1365

$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

$ws.Range("H38").Value = 1765.5294
$ws.Range("I38").Value = 139.76923
$ws.Range("K38").Value = 419.30769
$ws.Range("M38").Value = -47.30768999999998

$ws.Range("H107").Value = 2314.3333
$ws.Range("I107").Value = 1499.5
$ws.Range("J107").Value = 2721.75
$ws.Range("K107").Value = 1499.5
$ws.Range("L107").Value = 2721.75
$ws.Range("M107").Value = 420.5
$ws.Range("N107").Value = -6561.75

$ws.Range("H132").Value = 5094.448
$ws.Range("I132").Value = 2242.3
$ws.Range("J132").Value = 11432.556
$ws.Range("K132").Value = 6726.900000000001
$ws.Range("L132").Value = 34297.66800000001
$ws.Range("M132").Value = -4196.900000000001
$ws.Range("N132").Value = -39357.66800000001

$ws.Range("H135").Value = 1316.16
$ws.Range("I135").Value = 498.72726
$ws.Range("J135").Value = 7310.6665
$ws.Range("K135").Value = 4488.54534
$ws.Range("L135").Value = 65795.9985
$ws.Range("M135").Value = -1953.54534
$ws.Range("N135").Value = -70865.9985

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3990.22
$ws.Range("I32").Value = 3948.1875
$ws.Range("J32").Value = 4999
$ws.Range("K32").Value = 3948.1875
$ws.Range("L32").Value = 4999
$ws.Range("M32").Value = -3661.1875
$ws.Range("N32").Value = -5573

$ws.Range("H45").Value = 2581.9
$ws.Range("I45").Value = 1385.875
$ws.Range("J45").Value = 7366
$ws.Range("K45").Value = 1385.875
$ws.Range("L45").Value = 7366
$ws.Range("M45").Value = -1008.875
$ws.Range("N45").Value = -8120

$ws.Range("H97").Value = 881.1905
$ws.Range("I97").Value = 797.91895
$ws.Range("J97").Value = 1497.4
$ws.Range("K97").Value = 797.91895
$ws.Range("L97").Value = 1497.4
$ws.Range("M97").Value = -301.91895
$ws.Range("N97").Value = -2489.4

$ws.Range("H132").Value = 1788654
$ws.Range("I132").Value = 2844.5
$ws.Range("J132").Value = 16670400
$ws.Range("K132").Value = 8533.5
$ws.Range("L132").Value = 50011200
$ws.Range("M132").Value = -6003.5
$ws.Range("N132").Value = -50016260

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 299
$ws.Range("J7").Value = 726
$ws.Range("L7").Value = 726
$ws.Range("N7").Value = -952

$ws.Range("H22").Value = 632.1667
$ws.Range("I22").Value = 533
$ws.Range("K22").Value = 533
$ws.Range("M22").Value = -183

$ws.Range("H31").Value = 45922092
$ws.Range("I31").Value = 55559864
$ws.Range("J31").Value = 2552134.5
$ws.Range("K31").Value = 55559864
$ws.Range("L31").Value = 2552134.5
$ws.Range("M31").Value = -55559569
$ws.Range("N31").Value = -2552724.5

$ws.Range("H34").Value = 45922092
$ws.Range("I34").Value = 55559864
$ws.Range("J34").Value = 2552134.5
$ws.Range("K34").Value = 55559864
$ws.Range("L34").Value = 2552134.5
$ws.Range("M34").Value = -55559662
$ws.Range("N34").Value = -2552538.5

$ws.Range("H58").Value = 2249.3
$ws.Range("I58").Value = 2183.3333
$ws.Range("K58").Value = 2183.3333
$ws.Range("M58").Value = -1980.3333

$ws.Range("H62").Value = 4999.5
$ws.Range("J62").Value = 4999.5
$ws.Range("L62").Value = 4999.5
$ws.Range("N62").Value = -6247.5

$ws.Range("H65").Value = 4999.5
$ws.Range("J65").Value = 4999.5
$ws.Range("L65").Value = 24997.5
$ws.Range("N65").Value = -31237.5

$ws.Range("H107").Value = 2923.5757
$ws.Range("I107").Value = 2595.36
$ws.Range("J107").Value = 3949.25
$ws.Range("K107").Value = 2595.36
$ws.Range("L107").Value = 3949.25
$ws.Range("M107").Value = -675.3600000000001
$ws.Range("N107").Value = -7789.25

$ws.Range("H122").Value = 4283.3335
$ws.Range("I122").Value = 4140
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 12420
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -9970
$ws.Range("N122").Value = -19900

$ws.Range("H132").Value = 3207.1667
$ws.Range("I132").Value = 2832.111
$ws.Range("K132").Value = 8496.332999999999
$ws.Range("M132").Value = -5966.332999999999

$ws.Range("H134").Value = 1957.1052
$ws.Range("I134").Value = 2023.8
$ws.Range("J134").Value = 1707
$ws.Range("K134").Value = 6071.4
$ws.Range("L134").Value = 5121
$ws.Range("M134").Value = -3536.4
$ws.Range("N134").Value = -10191

$ws.Range("H136").Value = 2249.3
$ws.Range("I136").Value = 2183.3333
$ws.Range("K136").Value = 6549.999899999999
$ws.Range("M136").Value = -3999.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1474.8182
$ws.Range("I23").Value = 99
$ws.Range("J23").Value = 1612.4
$ws.Range("K23").Value = 297
$ws.Range("L23").Value = 4837.200000000001
$ws.Range("M23").Value = -62
$ws.Range("N23").Value = -5307.200000000001

$ws.Range("H26").Value = 587.8
$ws.Range("J26").Value = 495
$ws.Range("L26").Value = 1485
$ws.Range("N26").Value = -2061

$ws.Range("H33").Value = 5330918
$ws.Range("I33").Value = 187.25
$ws.Range("K33").Value = 1123.5
$ws.Range("M33").Value = -840.5

$ws.Range("H36").Value = 9254.333
$ws.Range("I36").Value = 4438.6
$ws.Range("K36").Value = 13315.8
$ws.Range("M36").Value = -13146.8

$ws.Range("H38").Value = 701.4

$ws.Range("H107").Value = 4140331.5
$ws.Range("I107").Value = 3815.6
$ws.Range("J107").Value = 5356954
$ws.Range("K107").Value = 11446.8
$ws.Range("L107").Value = 16070862
$ws.Range("M107").Value = -9526.8
$ws.Range("N107").Value = -16074702

$ws.Range("H141").Value = 8002.1113
$ws.Range("I141").Value = 4835.75
$ws.Range("K141").Value = 14507.25
$ws.Range("M141").Value = -9327.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 2619.8
$ws.Range("I107").Value = 2700
$ws.Range("K107").Value = 2700
$ws.Range("M107").Value = -780

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 99999
$ws.Range("J38").Value = 99999
$ws.Range("L38").Value = 99999
$ws.Range("N38").Value = -100819

$ws.Range("H39").Value = 39999
$ws.Range("J39").Value = 39999
$ws.Range("L39").Value = 39999
$ws.Range("N39").Value = -40919

$ws.Range("H40").Value = 6825.75
$ws.Range("J40").Value = 6500
$ws.Range("L40").Value = 6500
$ws.Range("N40").Value = -6772

$ws.Range("H136").Value = 3493.6365
$ws.Range("I136").Value = 2842.6924
$ws.Range("K136").Value = 8528.0772
$ws.Range("M136").Value = -5978.0772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 17512
$ws.Range("J101").Value = 17512
$ws.Range("L101").Value = 17512
$ws.Range("N101").Value = -24002

$ws.Range("H112").Value = 40899.668
$ws.Range("J112").Value = 40899.668
$ws.Range("L112").Value = 40899.668
$ws.Range("N112").Value = -43853.668

$ws.Range("H132").Value = 626962.06
$ws.Range("I132").Value = 2030.3077
$ws.Range("K132").Value = 6090.9231
$ws.Range("M132").Value = -3560.9231

$ws.Range("H136").Value = 279090.4
$ws.Range("I136").Value = 1373.5294
$ws.Range("J136").Value = 5000277.5
$ws.Range("K136").Value = 4120.5882
$ws.Range("L136").Value = 15000832.5
$ws.Range("M136").Value = -1570.5882
$ws.Range("N136").Value = -15005932.5
